# Add results and unfolding with 100 keV threshold
# Updates the Mean/HDI statistics for the existing "Beta" (row 2) and
# "Gamma" (row 3) results, and appends a new "Beta + Gamma" combined
# results row (row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ("Beta") updated values ----
$ws.Range("C2").Value = 40.35273080373261
$ws.Range("E2").Value = 0.05768901818751787
$ws.Range("F2").Value = 34.6928741187416
$ws.Range("G2").Value = 32.34012333811126
$ws.Range("H2").Value = 37.12535548611606
$ws.Range("I2").Value = 0.004882757102988988
$ws.Range("J2").Value = 0.0007258133397975762
$ws.Range("K2").Value = 0.01319101658353342
$ws.Range("L2").Value = 0.05697857888526592
$ws.Range("M2").Value = 0.0546420270339088
$ws.Range("N2").Value = 0.0599160276657242

# ---- Row 3 ("Gamma") updated values ----
$ws.Range("F3").Value = 0.5260965591264744
$ws.Range("G3").Value = 0.02687817300514678
$ws.Range("H3").Value = 1.044092587235455
$ws.Range("I3").Value = 0.4912437483309917
$ws.Range("J3").Value = 0.02482324644628219
$ws.Range("K3").Value = 0.974440389061549
$ws.Range("L3").Value = 0.5461418249637824
$ws.Range("M3").Value = 0.02830073138671181
$ws.Range("N3").Value = 1.083435470667743

# ---- New Row 4 ("Beta + Gamma") ----
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 40.35273080373261
$ws.Range("D4").Value = 0.0007432820064133916
$ws.Range("E4").Value = 0.05768901818751787
$ws.Range("F4").Value = 35.21897067786807
$ws.Range("G4").Value = 32.36700151111641
$ws.Range("H4").Value = 38.16944807335152
$ws.Range("I4").Value = 0.4961265054339807
$ws.Range("J4").Value = 0.02554905978607976
$ws.Range("K4").Value = 0.9876314056450823
$ws.Range("L4").Value = 0.6031204038490483
$ws.Range("M4").Value = 0.08294275842062061
$ws.Range("N4").Value = 1.143351498333467

# Match the bold/border/centered formatting used on A2 and A3 for the new A4 cell
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
